# "worked design review ppt"
#
# Travis logged a new time-card entry (9/8/2017, 3:00-4:00 PM, at "Rayzor Hall",
# working on "Prepared Design Review #1 Presentation"). The entry is inserted as
# a new row 8 on the "Travis" sheet, pushing the existing "Week 2 Total" row down
# to row 9 (its SUM range stays C5:C6, unaffected by the new row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Travis")

# Insert a fresh row above the current "Week 2 Total" row (row 8), shifting it to row 9.
$ws.Rows.Item(8).Insert()

# Seed the new row's formatting from the row above (the last real time entry, row 7)
# so the new cells pick up the same date/number formats, fonts, borders and wrap
# settings used throughout the sheet instead of generic defaults.
$ws.Range("A7:F7").Copy()
$ws.Range("A8:F8").PasteSpecial(-4122)

# New time-card entry values/formula.
$ws.Range("A8").Value = 42986.625
$ws.Range("B8").Value = 42986.666666666664
$ws.Range("C8").Formula = "=B8-A8"
$ws.Range("D8").Value = "Rayzor Hall"
$ws.Range("F8").Value = "Prepared Design Review #1 Presentation"

# The "Where?" gap column (E) on data rows lines up with the borderless style used
# in the Week Total rows; match it, and drop the border on the new "What did you
# work on?" cell (F8) to mirror how that column looks on the totals rows.
$ws.Range("E9").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("F8").Borders.LineStyle = 0

# Leave the cursor where the author left off after typing the new note.
[void]$ws.Range("F12").Select()
